$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking values are stored as text (matching original inlineStr type)
$textCells = @("D4", "D6", "D7", "D8", "D9", "D10", "D12", "D13", "D14", "D15", "D17", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D31", "D34", "D35", "D36", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply cell value updates per diff
$ws.Range("D2").Value = "95.721.84"
$ws.Range("E2").Value = "  -2.02%  "
$ws.Range("D3").Value = "3.623.52"
$ws.Range("E3").Value = "  -2.17%  "
$ws.Range("D4").Value = "2.72"
$ws.Range("E4").Value = "  +25.18%  "
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").Value = "223.52"
$ws.Range("E6").Value = "  -5.87%  "
$ws.Range("D7").Value = "642.94"
$ws.Range("E7").Value = "  -2.09%  "
$ws.Range("D8").Value = "0.424"
$ws.Range("E8").Value = "  -5.36%  "
$ws.Range("D9").Value = "1.20"
$ws.Range("E9").Value = "  +5.96%  "
$ws.Range("D10").Value = "1.00"
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("D11").Value = "3.621.39"
$ws.Range("E11").Value = "  -2.17%  "
$ws.Range("D12").Value = "51.55"
$ws.Range("E12").Value = "  +15.24%  "
$ws.Range("D13").Value = "0.218"
$ws.Range("E13").Value = "  +4.85%  "
$ws.Range("D14").Value = "0.0000292"
$ws.Range("E14").Value = "  -7.12%  "
$ws.Range("D15").Value = "6.54"
$ws.Range("E15").Value = "  -4.36%  "
$ws.Range("D16").Value = "4.295.79"
$ws.Range("E16").Value = "  -2.24%  "
$ws.Range("D17").Value = "24.89"
$ws.Range("E17").Value = "  +32.61%  "
$ws.Range("D18").Value = "95.510.35"
$ws.Range("E18").Value = "  -1.92%  "
$ws.Range("D19").Value = "9.23"
$ws.Range("E19").Value = "  +3.79%  "
$ws.Range("D20").Value = "13.85"
$ws.Range("E20").Value = "  +6.38%  "
$ws.Range("D21").Value = "3.616.09"
$ws.Range("E21").Value = "  -2.27%  "
$ws.Range("D22").Value = "0.536"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").Value = "0.282"
$ws.Range("E23").Value = "  +33.22%  "
$ws.Range("D24").Value = "137.38"
$ws.Range("E24").Value = "  +16.10%  "
$ws.Range("D25").Value = "536.05"
$ws.Range("E25").Value = "  +2.21%  "
$ws.Range("D26").Value = "3.28"
$ws.Range("E26").Value = "  -4.54%  "
$ws.Range("D27").Value = "0.0000203"
$ws.Range("E27").Value = "  -8.93%  "
$ws.Range("D28").Value = "6.97"
$ws.Range("E28").Value = "  +0.94%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "13.21"
$ws.Range("E29").Value = "  -1.90%  "
$ws.Range("B30").Value = "WrappedeETH"
$ws.Range("C30").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D30").Value = "3.788.75"
$ws.Range("E30").Value = "  -2.86%  "
$ws.Range("D31").Value = "13.44"
$ws.Range("E31").Value = "  +5.90%  "
$ws.Range("E32").Value = "  +3.98%  "
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "34.02"
$ws.Range("E34").Value = "  +3.31%  "
$ws.Range("B35").Value = "PolygonEcosystemToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D35").Value = "0.640"
$ws.Range("E35").Value = "  +7.53%  "
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").Value = "1.87"
$ws.Range("E36").Value = "  +2.76%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").Value = "7.34"
$ws.Range("E39").Value = "  +7.63%  "
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").Value = "8.54"
$ws.Range("E41").Value = "  -2.20%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "598.52"
$ws.Range("E42").Value = "  -6.06%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.0544"
$ws.Range("E43").Value = "  +20.06%  "
$ws.Range("E44").Value = "  +8.11%  "
$ws.Range("D45").Value = "0.506"
$ws.Range("E45").Value = "  +1.63%  "
$ws.Range("D46").Value = "41.31"
$ws.Range("E46").Value = "  +3.26%  "
$ws.Range("D47").Value = "0.162"
$ws.Range("E47").Value = "  -3.00%  "
$ws.Range("D48").Value = "2.00"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").Value = "9.46"
$ws.Range("E49").Value = "  +7.55%  "
$ws.Range("D50").Value = "232.38"
$ws.Range("E50").Value = "  +13.14%  "
$ws.Range("E51").Value = "  -1.87%  "

# Restore default (Normal) style on previously-reformatted cells to avoid leftover text-format styling
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
